$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new ranking entry at row 22 (Rafael Andrés / MDRplayer) ---
$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = "625190edabc36f0039acb8b2"
$ws.Range("B22").Value = "Rafael Andrés"
$ws.Range("C22").Value = 353140182.88849998
$ws.Range("D22").Value = 21
$ws.Range("E22").Value = "MDRplayer"

# --- Insert new ranking entry at row 38 (Leandro / lw2idgxo) ---
$ws.Rows.Item(38).Insert()
$ws.Range("A38").Value = "663fc91267e60ee7d8a2f703"
$ws.Range("B38").Value = "Leandro"
$ws.Range("C38").Value = 85943046.244000003
$ws.Range("D38").Value = 37
$ws.Range("E38").Value = "lw2idgxo"

# --- Fix up the "posicao" (D) column sequence for every row pushed down ---
# Row 23..37 held the old rows 22..36 (posicao 21..35) -> need 22..36
for ($r = 23; $r -le 37; $r++) {
    $ws.Cells.Item($r, 4).Value = $r - 1
}
# Row 39..51 held the old rows 37..49 (posicao 36..48) -> need 38..50
for ($r = 39; $r -le 51; $r++) {
    $ws.Cells.Item($r, 4).Value = $r - 1
}

# --- Update the view state to match the saved workbook ---
$ws.Range("A2:E51").Select()
$excel.ActiveWindow.ScrollRow = 22

Write-Output "done"
